$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.03877914339252267
    "C2" = 0.1223661001864929
    "D2" = 0.5921565300001895
    "E2" = 0.2580542788517324
    "F2" = 0.09150221478225759
    "G2" = 0.0475100079502194
    "H2" = 0.1099094607512485
    "B3" = 0.1894635701346697
    "C3" = 0.6592539999483663
    "D3" = 0.3251517487999092
    "E3" = 0.1585996847304344
    "F3" = 0.1146074778983962
    "G3" = 0.1770069306994253
    "B4" = 0.5408953042689693
    "C4" = 0.2067930531205122
    "D4" = 0.04024098905103731
    "E4" = -0.00375121778100088
    "F4" = 0.05864823502002824
    "G4" = 0.02541975869265006
    "H4" = 0.03155836416756893
    "I4" = -0.02336658152311553
    "J4" = -0.02233679963354765
    "B5" = 0.1798140101391425
    "C5" = 0.01326194606966766
    "D5" = -0.03073026076237053
    "E5" = 0.03166919203865859
    "F5" = -0.001559284288719586
    "G5" = 0.004579321186199281
    "H5" = -0.05034562450448518
    "I5" = -0.0493158426149173
    "B6" = 0.2866721711845134
    "C6" = 0.2426799643524752
    "D6" = 0.3050794171535043
    "E6" = 0.2718509408261262
    "F6" = 0.277989546301045
    "G6" = 0.2230646006103605
    "H6" = 0.2240943824999284
    "B7" = 0.2293351707594228
    "C7" = 0.2917346235604519
    "D7" = 0.2585061472330737
    "E7" = 0.2646447527079926
    "F7" = 0.2097198070173081
    "G7" = 0.210749588906876
    "B8" = 0.03613384424365556
    "C8" = 0.002905367916277386
    "D8" = 0.009043973391196253
    "E8" = -0.04588097229948821
    "F8" = -0.04485119040992033
    "G8" = -0.4685329518921513
    "H8" = 0.02897294484477524
    "I8" = -0.04257272378961845
    "B9" = 0.09666542263617448
    "C9" = 0.1028040281110933
    "D9" = 0.04787908242040887
    "E9" = 0.04890886430997676
    "F9" = -0.3747728971722541
    "G9" = 0.1227329995646723
    "H9" = 0.05118733093027864
    "B10" = 0.05612687574383117
    "C10" = 0.001201930053146706
    "D10" = 0.002231711942714591
    "E10" = -0.4214500495395163
    "F10" = 0.07605584719741015
    "G10" = 0.00451017856301647
    "B11" = -0.06941193049434516
    "C11" = -0.06838214860477726
    "D11" = -0.4920639100870082
    "E11" = 0.005441986649918303
    "F11" = -0.06610368198447539
    "B12" = -0.002889654035708851
    "C12" = -0.4265714155179398
    "D12" = 0.07093448121898671
    "E12" = -0.0006111874154069719
    "B13" = -0.4066743937068965
    "C13" = 0.09083150303002996
    "D13" = 0.01928583439563626
    "B14" = 0.1531996791782531
    "C14" = 0.08165401054385939
    "B15" = -0.1180965791298333
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
